# The edit rotates the species-occurrence records stored in sheet rows 4-10:
# each row's identifying data moves down into the next row, and the record
# that used to be on the last row (10) wraps around to become the new first
# row (4). Concretely (old Id -> new row):
#   old A10 (102089658) -> new A4
#   old A4  (102089519) -> new A5
#   old A5  (102089499) -> new A6
#   old A6  (102089518) -> new A7
#   old A7  (102089652) -> new A8
#   old A8  (102089546) -> new A9
#   old A9  (102089632) -> new A10
#
# Only columns A, B, D, E, F, G, H, Q and R actually hold values that differ
# from row to row in this range (everything else - C, I, P, S, T, U, V, W,
# Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY, etc. - is identical across all
# seven rows), so only those columns need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colsToRotate = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $colsToRotate) {
    $rng = $ws.Range("${col}4:${col}10")
    $val = $rng.Value()
    $rowCount = $val.GetLength(0)

    # Range.Value() returns a 1-based array, but assigning back to
    # Range.Value expects a 0-based array, so build the replacement
    # accordingly.
    $newVal = New-Object 'object[,]' $rowCount, 1

    # New first row (sheet row 4) = old last row (sheet row 10)
    $newVal[0, 0] = $val[$rowCount, 1]

    # Every remaining row shifts down by one: new row r = old row r-1
    for ($r = 2; $r -le $rowCount; $r++) {
        $newVal[$r - 1, 0] = $val[$r - 1, 1]
    }

    $rng.Value = $newVal
}
